$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value = $text
    $cellRange.Style = "Normal"
}

$ws.Range("D2").Value = '60.255.24'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.411.06'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  +0.00%  '
Set-TextValue $ws.Range("D5") '559.37'
$ws.Range("E5").Value = '  +1.39%  '
Set-TextValue $ws.Range("D6") '135.74'
$ws.Range("E6").Value = '  -1.46%  '
$ws.Range("E7").Value = '  +0.01%  '
Set-TextValue $ws.Range("D8") '0.589'
$ws.Range("E8").Value = '  +0.41%  '
$ws.Range("E9").Value = '  +0.28%  '
$ws.Range("E10").Value = '  -1.27%  '
$ws.Range("E11").Value = '  +1.23%  '
$ws.Range("E12").Value = '  -1.62%  '
Set-TextValue $ws.Range("D13") '24.70'
$ws.Range("E13").Value = '  -3.44%  '
$ws.Range("D14").Value = '2.837.44'
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("D15").Value = '60.160.37'
$ws.Range("E15").Value = '  +0.50%  '
$ws.Range("E16").Value = '  +0.51%  '
$ws.Range("D17").Value = '2.374.28'
$ws.Range("E17").Value = '  -2.39%  '
Set-TextValue $ws.Range("D18") '11.22'
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("E19").Value = '  +3.25%  '
Set-TextValue $ws.Range("D20") '325.04'
$ws.Range("E20").Value = '  -1.83%  '
$ws.Range("E21").Value = '  +1.71%  '
$ws.Range("E22").Value = '  +0.09%  '
Set-TextValue $ws.Range("D23") '64.62'
$ws.Range("E23").Value = '  -2.58%  '
$ws.Range("E24").Value = '  +1.48%  '
Set-TextValue $ws.Range("D25") '8.53'
$ws.Range("E25").Value = '  -2.86%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").Value = '0.0₃0774'
$ws.Range("E29").Value = '  -0.77%  '
Set-TextValue $ws.Range("D30") '170.68'
$ws.Range("E30").Value = '  +1.20%  '
Set-TextValue $ws.Range("D31") '6.13'
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("E32").Value = '  +5.80%  '
$ws.Range("E33").Value = '  -2.35%  '
Set-TextValue $ws.Range("D34") '18.37'
$ws.Range("E34").Value = '  -1.70%  '
$ws.Range("E35").Value = '  +3.89%  '
$ws.Range("E36").Value = '  +0.00%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  -1.22%  '
Set-TextValue $ws.Range("D39") '324.11'
$ws.Range("E39").Value = '  +3.10%  '
$ws.Range("E40").Value = '  -0.60%  '
Set-TextValue $ws.Range("D41") '38.64'
$ws.Range("E41").Value = '  -2.51%  '
Set-TextValue $ws.Range("D42") '147.92'
$ws.Range("E42").Value = '  +6.17%  '
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("E44").Value = '  +0.17%  '
Set-TextValue $ws.Range("D45") '19.90'
$ws.Range("E45").Value = '  +2.10%  '
Set-TextValue $ws.Range("D46") '0.0517'
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("E47").Value = '  -0.58%  '
Set-TextValue $ws.Range("D48") '0.0222'
$ws.Range("E48").Value = '  -1.50%  '
$ws.Range("E49").Value = '  -0.16%  '
Set-TextValue $ws.Range("D50") '1.58'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("E51").Value = '  -0.61%  '
